$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).NumberFormat = '@'
$ws.Cells.Item(2, 4).Value = '23.457.18'
$ws.Cells.Item(2, 4).Style = 'Normal'
$ws.Cells.Item(2, 5).Value = '  -0.26%  '

# Row 3
$ws.Cells.Item(3, 4).NumberFormat = '@'
$ws.Cells.Item(3, 4).Value = '1.647.94'
$ws.Cells.Item(3, 4).Style = 'Normal'
$ws.Cells.Item(3, 5).Value = '  +0.13%  '

# Row 4
$ws.Cells.Item(4, 4).NumberFormat = '@'
$ws.Cells.Item(4, 4).Value = '0.9998'
$ws.Cells.Item(4, 4).Style = 'Normal'
$ws.Cells.Item(4, 5).Value = '  -0.39%  '

# Row 5
$ws.Cells.Item(5, 4).NumberFormat = '@'
$ws.Cells.Item(5, 4).Value = '0.9995'
$ws.Cells.Item(5, 4).Style = 'Normal'
$ws.Cells.Item(5, 5).Value = '  -0.41%  '

# Row 6
$ws.Cells.Item(6, 4).NumberFormat = '@'
$ws.Cells.Item(6, 4).Value = '299.94'
$ws.Cells.Item(6, 4).Style = 'Normal'
$ws.Cells.Item(6, 5).Value = '  -0.87%  '

# Row 7
$ws.Cells.Item(7, 4).NumberFormat = '@'
$ws.Cells.Item(7, 4).Value = '0.3797'
$ws.Cells.Item(7, 4).Style = 'Normal'
$ws.Cells.Item(7, 5).Value = '  -1.13%  '

# Row 8
$ws.Cells.Item(8, 4).NumberFormat = '@'
$ws.Cells.Item(8, 4).Value = '50.48'
$ws.Cells.Item(8, 4).Style = 'Normal'
$ws.Cells.Item(8, 5).Value = '  -0.94%  '

# Row 9
$ws.Cells.Item(9, 4).NumberFormat = '@'
$ws.Cells.Item(9, 4).Value = '0.3504'
$ws.Cells.Item(9, 4).Style = 'Normal'
$ws.Cells.Item(9, 5).Value = '  -2.57%  '

# Row 10
$ws.Cells.Item(10, 4).NumberFormat = '@'
$ws.Cells.Item(10, 4).Value = '1.222'
$ws.Cells.Item(10, 4).Style = 'Normal'
$ws.Cells.Item(10, 5).Value = '  -0.44%  '

# Row 11
$ws.Cells.Item(11, 4).NumberFormat = '@'
$ws.Cells.Item(11, 4).Value = '0.08074'
$ws.Cells.Item(11, 4).Style = 'Normal'
$ws.Cells.Item(11, 5).Value = '  -1.12%  '

# Row 12
$ws.Cells.Item(12, 4).NumberFormat = '@'
$ws.Cells.Item(12, 4).Value = '0.9998'
$ws.Cells.Item(12, 4).Style = 'Normal'
$ws.Cells.Item(12, 5).Value = '  -0.42%  '

# Row 13
$ws.Cells.Item(13, 4).NumberFormat = '@'
$ws.Cells.Item(13, 4).Value = '22.07'
$ws.Cells.Item(13, 4).Style = 'Normal'
$ws.Cells.Item(13, 5).Value = '  -0.97%  '

# Row 14
$ws.Cells.Item(14, 4).NumberFormat = '@'
$ws.Cells.Item(14, 4).Value = '6.321'
$ws.Cells.Item(14, 4).Style = 'Normal'
$ws.Cells.Item(14, 5).Value = '  -1.71%  '

# Row 15
$ws.Cells.Item(15, 4).NumberFormat = '@'
$ws.Cells.Item(15, 4).Value = '7.265'
$ws.Cells.Item(15, 4).Style = 'Normal'
$ws.Cells.Item(15, 5).Value = '  -2.21%  '

# Row 16
$ws.Cells.Item(16, 4).NumberFormat = '@'
$ws.Cells.Item(16, 4).Value = '0.00001214'
$ws.Cells.Item(16, 4).Style = 'Normal'
$ws.Cells.Item(16, 5).Value = '  -0.25%  '

# Row 17
$ws.Cells.Item(17, 4).NumberFormat = '@'
$ws.Cells.Item(17, 4).Value = '1.662.57'
$ws.Cells.Item(17, 4).Style = 'Normal'
$ws.Cells.Item(17, 5).Value = '  +0.69%  '

# Row 18
$ws.Cells.Item(18, 4).NumberFormat = '@'
$ws.Cells.Item(18, 4).Value = '94.91'
$ws.Cells.Item(18, 4).Style = 'Normal'
$ws.Cells.Item(18, 5).Value = '  -2.58%  '

# Row 19
$ws.Cells.Item(19, 4).NumberFormat = '@'
$ws.Cells.Item(19, 4).Value = '0.06972'
$ws.Cells.Item(19, 4).Style = 'Normal'
$ws.Cells.Item(19, 5).Value = '  -0.60%  '

# Row 20
$ws.Cells.Item(20, 4).NumberFormat = '@'
$ws.Cells.Item(20, 4).Value = '6.626'
$ws.Cells.Item(20, 4).Style = 'Normal'
$ws.Cells.Item(20, 5).Value = '  -1.89%  '

# Row 21
$ws.Cells.Item(21, 4).NumberFormat = '@'
$ws.Cells.Item(21, 4).Value = '17.45'
$ws.Cells.Item(21, 4).Style = 'Normal'
$ws.Cells.Item(21, 5).Value = '  -0.36%  '

# Row 22
$ws.Cells.Item(22, 4).NumberFormat = '@'
$ws.Cells.Item(22, 4).Value = '0.9992'
$ws.Cells.Item(22, 4).Style = 'Normal'
$ws.Cells.Item(22, 5).Value = '  -0.44%  '

# Row 23
$ws.Cells.Item(23, 5).Value = '  -1.15%  '

# Row 24
$ws.Cells.Item(24, 4).NumberFormat = '@'
$ws.Cells.Item(24, 4).Value = '23.462.56'
$ws.Cells.Item(24, 4).Style = 'Normal'
$ws.Cells.Item(24, 5).Value = '  -0.26%  '

# Row 25
$ws.Cells.Item(25, 4).NumberFormat = '@'
$ws.Cells.Item(25, 4).Value = '2.427'
$ws.Cells.Item(25, 4).Style = 'Normal'
$ws.Cells.Item(25, 5).Value = '  -2.10%  '

# Row 26
$ws.Cells.Item(26, 4).NumberFormat = '@'
$ws.Cells.Item(26, 4).Value = '2.989'
$ws.Cells.Item(26, 4).Style = 'Normal'
$ws.Cells.Item(26, 5).Value = '  -1.31%  '

# Row 27
$ws.Cells.Item(27, 4).NumberFormat = '@'
$ws.Cells.Item(27, 4).Value = '21.08'
$ws.Cells.Item(27, 4).Style = 'Normal'
$ws.Cells.Item(27, 5).Value = '  -0.45%  '

# Row 28
$ws.Cells.Item(28, 4).NumberFormat = '@'
$ws.Cells.Item(28, 4).Value = '151.41'
$ws.Cells.Item(28, 4).Style = 'Normal'
$ws.Cells.Item(28, 5).Value = '  -1.38%  '

# Row 29
$ws.Cells.Item(29, 4).NumberFormat = '@'
$ws.Cells.Item(29, 4).Value = '5.212'
$ws.Cells.Item(29, 4).Style = 'Normal'
$ws.Cells.Item(29, 5).Value = '  -0.36%  '

# Row 30
$ws.Cells.Item(30, 4).NumberFormat = '@'
$ws.Cells.Item(30, 4).Value = '131.93'
$ws.Cells.Item(30, 4).Style = 'Normal'
$ws.Cells.Item(30, 5).Value = '  -1.37%  '

# Row 31
$ws.Cells.Item(31, 4).NumberFormat = '@'
$ws.Cells.Item(31, 4).Value = '1.831.13'
$ws.Cells.Item(31, 4).Style = 'Normal'
$ws.Cells.Item(31, 5).Value = '  -0.04%  '

# Row 32
$ws.Cells.Item(32, 4).NumberFormat = '@'
$ws.Cells.Item(32, 4).Value = '6.875'
$ws.Cells.Item(32, 4).Style = 'Normal'
$ws.Cells.Item(32, 5).Value = '  -2.35%  '

# Row 33
$ws.Cells.Item(33, 4).NumberFormat = '@'
$ws.Cells.Item(33, 4).Value = '2.135'
$ws.Cells.Item(33, 4).Style = 'Normal'
$ws.Cells.Item(33, 5).Value = '  -5.11%  '

# Row 34
$ws.Cells.Item(34, 4).NumberFormat = '@'
$ws.Cells.Item(34, 4).Value = '11.17'
$ws.Cells.Item(34, 4).Style = 'Normal'
$ws.Cells.Item(34, 5).Value = '  -8.27%  '

# Row 35
$ws.Cells.Item(35, 4).NumberFormat = '@'
$ws.Cells.Item(35, 4).Value = '0.9897'
$ws.Cells.Item(35, 4).Style = 'Normal'
$ws.Cells.Item(35, 5).Value = '  -6.08%  '

# Row 36
$ws.Cells.Item(36, 4).NumberFormat = '@'
$ws.Cells.Item(36, 4).Value = '0.02701'
$ws.Cells.Item(36, 4).Style = 'Normal'
$ws.Cells.Item(36, 5).Value = '  -3.02%  '

# Row 37
$ws.Cells.Item(37, 4).NumberFormat = '@'
$ws.Cells.Item(37, 4).Value = '0.08792'
$ws.Cells.Item(37, 4).Style = 'Normal'
$ws.Cells.Item(37, 5).Value = '  +0.24%  '

# Row 38
$ws.Cells.Item(38, 2).Value = 'InternetComputer(DFINITY)'
$ws.Cells.Item(38, 3).Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Cells.Item(38, 4).NumberFormat = '@'
$ws.Cells.Item(38, 4).Value = '5.919'
$ws.Cells.Item(38, 4).Style = 'Normal'
$ws.Cells.Item(38, 5).Value = '  -2.19%  '

# Row 39
$ws.Cells.Item(39, 2).Value = 'Algorand'
$ws.Cells.Item(39, 3).Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Cells.Item(39, 4).NumberFormat = '@'
$ws.Cells.Item(39, 4).Value = '0.2423'
$ws.Cells.Item(39, 4).Style = 'Normal'
$ws.Cells.Item(39, 5).Value = '  -2.91%  '

# Row 40
$ws.Cells.Item(40, 4).NumberFormat = '@'
$ws.Cells.Item(40, 4).Value = '0.06806'
$ws.Cells.Item(40, 4).Style = 'Normal'
$ws.Cells.Item(40, 5).Value = '  -2.29%  '

# Row 41
$ws.Cells.Item(41, 4).NumberFormat = '@'
$ws.Cells.Item(41, 4).Value = '12.80'
$ws.Cells.Item(41, 4).Style = 'Normal'
$ws.Cells.Item(41, 5).Value = '  -2.10%  '

# Row 42
$ws.Cells.Item(42, 4).NumberFormat = '@'
$ws.Cells.Item(42, 4).Value = '0.6857'
$ws.Cells.Item(42, 4).Style = 'Normal'
$ws.Cells.Item(42, 5).Value = '  -1.57%  '

# Row 43
$ws.Cells.Item(43, 5).Value = '  -3.04%  '

# Row 44
$ws.Cells.Item(44, 5).Value = '  -2.51%  '

# Row 45
$ws.Cells.Item(45, 4).NumberFormat = '@'
$ws.Cells.Item(45, 4).Value = '0.9985'
$ws.Cells.Item(45, 4).Style = 'Normal'
$ws.Cells.Item(45, 5).Value = '  -0.45%  '

# Row 46
$ws.Cells.Item(46, 4).NumberFormat = '@'
$ws.Cells.Item(46, 4).Value = '0.6375'
$ws.Cells.Item(46, 4).Style = 'Normal'
$ws.Cells.Item(46, 5).Value = '  -1.57%  '

# Row 47
$ws.Cells.Item(47, 4).NumberFormat = '@'
$ws.Cells.Item(47, 4).Value = '2.255'
$ws.Cells.Item(47, 4).Style = 'Normal'
$ws.Cells.Item(47, 5).Value = '  -1.44%  '

# Row 48
$ws.Cells.Item(48, 4).NumberFormat = '@'
$ws.Cells.Item(48, 4).Value = '3.914'
$ws.Cells.Item(48, 4).Style = 'Normal'
$ws.Cells.Item(48, 5).Value = '  -1.09%  '

# Row 49
$ws.Cells.Item(49, 4).NumberFormat = '@'
$ws.Cells.Item(49, 4).Value = '0.07684'
$ws.Cells.Item(49, 4).Style = 'Normal'

# Row 50
$ws.Cells.Item(50, 4).NumberFormat = '@'
$ws.Cells.Item(50, 4).Value = '126.90'
$ws.Cells.Item(50, 4).Style = 'Normal'
$ws.Cells.Item(50, 5).Value = '  -0.71%  '

# Row 51
$ws.Cells.Item(51, 4).NumberFormat = '@'
$ws.Cells.Item(51, 4).Value = '1.230'
$ws.Cells.Item(51, 4).Style = 'Normal'
